$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (header): collapse Min/Mid/Max/PrefaultTime columns into one TriggerPriority column ---
$ws.Range("G1").Value = "'TriggerPriority"
$ws.Range("H1").Value = ""
$ws.Range("I1").Value = ""
$ws.Range("J1").Value = ""

# --- Row 2: first boundary-value row (TriggerPriority = 1) ---
$ws.Range("G2").Value = "'1"
$ws.Range("H2").Value = ""
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = ""

# --- Rows 3-6: new boundary-value test rows, duplicating the device info from row 2 ---
# Use a full-range Copy (preserves text styles + values in one shot), then patch column F's
# style back with a formats-only paste (a bare numeric .Value assignment resets cell style).
$ws.Range("A2:K2").Copy($ws.Range("A3:K3"))
$ws.Range("A2:K2").Copy($ws.Range("A4:K4"))
$ws.Range("A2:K2").Copy($ws.Range("A5:K5"))
$ws.Range("A2:K2").Copy($ws.Range("A6:K6"))

$ws.Range("F2").Copy()
$ws.Range("F3").PasteSpecial(-4122)
$ws.Range("F2").Copy()
$ws.Range("F4").PasteSpecial(-4122)
$ws.Range("F2").Copy()
$ws.Range("F5").PasteSpecial(-4122)
$ws.Range("F2").Copy()
$ws.Range("F6").PasteSpecial(-4122)

$ws.Range("G3").Value = "'99"
$ws.Range("G4").Value = "'15"
$ws.Range("G5").Value = "'24"
$ws.Range("G6").Value = "'64"

$ws.Range("G1").Select()
